$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 10:57"

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 795755
$ws.Range("C6").Value = 913
$ws.Range("D6").Value = 496048
$ws.Range("E6").Value = 278068
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 16
$ws.Range("H6").Value = 21639

# Row 7: Rusia
$ws.Range("A7").Value = "Rusia"
$ws.Range("B7").Value = 713936
$ws.Range("C7").Value = 6635
$ws.Range("D7").Value = 489068
$ws.Range("E7").Value = 213851
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 174
$ws.Range("H7").Value = 11017

# Row 20: Banglades
$ws.Range("A20").Value = "Banglades"
$ws.Range("B20").Value = 178443
$ws.Range("C20").Value = 2949
$ws.Range("D20").Value = 86406
$ws.Range("E20").Value = 89762
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = 2275

# Row 46: Polonia
$ws.Range("A46").Value = "Polonia"
$ws.Range("B46").Value = 37216
$ws.Range("C46").Value = 265
$ws.Range("D46").Value = 26048
$ws.Range("E46").Value = 9606
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 11
$ws.Range("H46").Value = 1562

# Row 47: Israel
$ws.Range("A47").Value = "Israel"
$ws.Range("B47").Value = 35631
$ws.Range("C47").Value = 806
$ws.Range("D47").Value = 18542
$ws.Range("E47").Value = 16739
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 350

# Row 60: Moldavia
$ws.Range("A60").Value = "Moldavia"
$ws.Range("B60").Value = 18666
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 12188
$ws.Range("E60").Value = 5851
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 627

# Row 75: El Salvador (now ranks above Kenia/Noruega)
$ws.Range("A75").Value = "El Salvador"
$ws.Range("B75").Value = 9142
$ws.Range("C75").Value = 298
$ws.Range("D75").Value = 5428
$ws.Range("E75").Value = 3465
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 6
$ws.Range("H75").Value = 249

# Row 76: Kenia (shifted down one row)
$ws.Range("A76").Value = "Kenia"
$ws.Range("B76").Value = 8975
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 2657
$ws.Range("E76").Value = 6145
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 173

# Row 77: Noruega (shifted down one row)
$ws.Range("A77").Value = "Noruega"
$ws.Range("B77").Value = 8965
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 8138
$ws.Range("E77").Value = 575
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 252

# Row 93: Estado de Palestina
$ws.Range("A93").Value = "Estado de Palestina"
$ws.Range("B93").Value = 5220
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 525
$ws.Range("E93").Value = 4668
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 27

# Row 119: Eslovaquia (now ranks above Lituania)
$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("B119").Value = 1870
$ws.Range("C119").Value = 19
$ws.Range("D119").Value = 1481
$ws.Range("E119").Value = 361
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 28

# Row 120: Lituania (shifted down one row)
$ws.Range("A120").Value = "Lituania"
$ws.Range("B120").Value = 1861
$ws.Range("C120").Value = 4
$ws.Range("D120").Value = 1569
$ws.Range("E120").Value = 213
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 79

# Row 122: Eslovenia (now ranks above Guinea-Bisau)
$ws.Range("A122").Value = "Eslovenia"
$ws.Range("B122").Value = 1793
$ws.Range("C122").Value = 17
$ws.Range("D122").Value = 1429
$ws.Range("E122").Value = 253
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 111

# Row 123: Guinea-Bisau (shifted down one row)
$ws.Range("A123").Value = "Guinea-Bisau"
$ws.Range("B123").Value = 1790
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 760
$ws.Range("E123").Value = 1005
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 25

# Row 135: Letonia
$ws.Range("A135").Value = "Letonia"
$ws.Range("B135").Value = 1165
$ws.Range("C135").Value = 11
$ws.Range("D135").Value = 1019
$ws.Range("E135").Value = 116
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 30

# Row 158: Taiwan
$ws.Range("A158").Value = "Taiwan"
$ws.Range("B158").Value = 451
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 438
$ws.Range("E158").Value = 6
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 7
